$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=-80.38; "C"=528.72; "D"=50; "E"=0.41; "F"=-82.82; "G"=528.72; "H"=50; "I"=0.4; "J"=-26.88; "K"=528.72; "L"=50; "M"=1.24; "N"=43392.34; "O"=113.4; "P"=51.23; "Q"=6.61 }
    3 = @{ "B"=-90.65; "C"=528.72; "D"=50; "E"=0.37; "F"=-89.84; "G"=528.72; "H"=50; "I"=0.37; "J"=-30.32; "K"=528.72; "L"=50; "M"=1.1; "N"=43392.34; "O"=113.4; "P"=51.23; "Q"=6.61 }
    4 = @{ "B"=-90.37; "C"=528.72; "D"=50; "E"=0.37; "F"=-90.11; "G"=528.72; "H"=50; "I"=0.37; "J"=-30.23; "K"=528.72; "L"=50; "M"=1.1; "N"=43392.34; "O"=113.4; "P"=51.23; "Q"=6.61 }
    5 = @{ "B"=-89.84; "C"=528.72; "D"=50; "E"=0.37; "F"=-89.86; "G"=528.72; "H"=50; "I"=0.37; "J"=-30.05; "K"=528.72; "L"=50; "M"=1.11; "N"=43392.34; "O"=113.4; "P"=51.23; "Q"=6.61 }
    6 = @{ "B"=-89.62; "C"=528.72; "D"=50; "E"=0.37; "F"=-89.68; "G"=528.72; "H"=50; "I"=0.37; "J"=-29.97; "K"=528.72; "L"=50; "M"=1.11; "N"=43392.34; "O"=113.4; "P"=51.23; "Q"=6.61 }
    7 = @{ "B"=-89.84; "C"=528.72; "D"=50; "E"=0.37; "F"=-89.86; "G"=528.72; "H"=50; "I"=0.37; "J"=-30.05; "K"=528.72; "L"=50; "M"=1.11; "N"=43392.34; "O"=113.4; "P"=51.23; "Q"=6.61 }
    8 = @{ "B"=-90.37; "C"=528.72; "D"=50; "E"=0.37; "F"=-90.11; "G"=528.72; "H"=50; "I"=0.37; "J"=-30.23; "K"=528.72; "L"=50; "M"=1.1; "N"=43392.34; "O"=113.4; "P"=51.23; "Q"=6.61 }
    9 = @{ "B"=-90.65; "C"=528.72; "D"=50; "E"=0.37; "F"=-89.84; "G"=528.72; "H"=50; "I"=0.37; "J"=-30.32; "K"=528.72; "L"=50; "M"=1.1; "N"=43392.34; "O"=113.4; "P"=51.23; "Q"=6.61 }
    10 = @{ "B"=-80.38; "C"=528.72; "D"=50; "E"=0.41; "F"=-82.82; "G"=528.72; "H"=50; "I"=0.4; "J"=-26.88; "K"=528.72; "L"=50; "M"=1.24; "N"=43392.34; "O"=113.4; "P"=51.23; "Q"=6.61 }
}

foreach ($r in $data.Keys) {
    foreach ($col in $data[$r].Keys) {
        $ws.Range("$col$r").Value = $data[$r][$col]
    }
}
